$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LOM3013")

# Ativação: 01/01/2015 -> 01/01/2023 (keep as literal text, not an auto-parsed date)
$ws.Range("B8:C8").NumberFormat = "@"
$ws.Range("B8").Value = "01/01/2023"
$ws.Range("C8").Value = "01/01/2023"

# Critério: short text -> long descriptive text
$criterio = "Esta é uma disciplina de caráter fundamental, exigindo dedicação individual para assimilação das definições e conceitos. Isto envolve leitura concentrada para fixação dos conceitos teóricos e realização de exercícios numéricos. Duas provas escritas (P1 e P2) serão aplicadas e com pesos iguais. O desenvolvimento do aluno ao longo do curso será aferido e estimulado por meio de discussões sobre um dado tema, porém sem a atribuição de nota, por conta da subjetividade envolvida."
$ws.Range("B18").Value = $criterio
$ws.Range("C18").Value = $criterio

# Norma de recuperação: update formula text
$norma = ": A Nota final (NF) será calculada da seguinte maneira: NF = (0,4*P1 +0,4* P2+ 0,2*NT) / 3"
$ws.Range("B19").Value = $norma
$ws.Range("C19").Value = $norma

# Bibliografia (recovery) text: add spacing before /2
$biblio = "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR) / 2"
$ws.Range("B20").Value = $biblio
$ws.Range("C20").Value = $biblio
